# "pretest and scan instructions update"
#
# 1) Workbook uses the 1904 date system (workbookPr date1904="1").
# 2) Row 8 ("In this experiment...") shrinks from 92.35 -> 56.35.
# 3) Two new instruction rows are inserted with text (using the row 8/9
#    formatting as a template so styles + row heights line up):
#      row 9  (56.35): "Even if we do not give you a perspective for
#                        some of them, please listen to all of the
#                        stories carefully. "
#      row 10 (92.35): "Try to be as still as possible during the scan,
#                        especially as you are listening to the
#                        stories. "
# 4) Two blank spacer rows (23, 24) are appended at the bottom, matching
#    the style/height of the existing blank rows (e.g. row 22).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- 1904 date system ------------------------------------------------
$wb.Date1904 = $true

# -- Row 8: new height -------------------------------------------------
$ws.Rows.Item(8).RowHeight = 56.35

# -- Row 9: copy row 8 formatting, then set its own text + height ------
$ws.Range("A8:G8").Copy()
$ws.Range("A9:G9").PasteSpecial(-4122)
$ws.Range("A9").Value = "Even if we do not give you a perspective for some of them, please listen to all of the stories carefully. "
$ws.Rows.Item(9).RowHeight = 56.35

# -- Row 10: copy row 8 formatting, then set its own text + height -----
$ws.Range("A8:G8").Copy()
$ws.Range("A10:G10").PasteSpecial(-4122)
$ws.Range("A10").Value = "Try to be as still as possible during the scan, especially as you are listening to the stories. "
$ws.Rows.Item(10).RowHeight = 92.35

# -- Rows 23 & 24: blank spacer rows, matching row 22's formatting -----
$ws.Range("A22:G22").Copy()
$ws.Range("A23:G23").PasteSpecial(-4122)
$ws.Rows.Item(23).RowHeight = 20.35

$ws.Range("A22:G22").Copy()
$ws.Range("A24:G24").PasteSpecial(-4122)
$ws.Rows.Item(24).RowHeight = 20.35
